$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 127.14286
$ws.Range("I2").Value = 127.14286
$ws.Range("K2").Value = 127.14286
$ws.Range("M2").Value = -14.14286
$ws.Range("H11").Value = 33.666668
$ws.Range("I11").Value = 33.666668
$ws.Range("K11").Value = 33.666668
$ws.Range("M11").Value = 106.333332
$ws.Range("H18").Value = 649.8333
$ws.Range("I18").Value = 379.8
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 379.8
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = -95.80000000000001
$ws.Range("N18").Value = -2568
$ws.Range("H55").Value = 668.25
$ws.Range("I55").Value = 558
$ws.Range("K55").Value = 558
$ws.Range("M55").Value = -344
$ws.Range("H70").Value = 460150.78
$ws.Range("J70").Value = 85035.836
$ws.Range("L70").Value = 255107.508
$ws.Range("N70").Value = -255647.508
$ws.Range("H73").Value = 460150.78
$ws.Range("J73").Value = 85035.836
$ws.Range("L73").Value = 255107.508
$ws.Range("N73").Value = -256979.508
$ws.Range("H92").Value = 323.7037
$ws.Range("I92").Value = 333.15
$ws.Range("J92").Value = 296.7143
$ws.Range("K92").Value = 333.15
$ws.Range("L92").Value = 296.7143
$ws.Range("M92").Value = 914.85
$ws.Range("N92").Value = -2792.7143
$ws.Range("H96").Value = 711.86664
$ws.Range("I96").Value = 668.36365
$ws.Range("J96").Value = 831.5
$ws.Range("K96").Value = 2005.09095
$ws.Range("L96").Value = 2494.5
$ws.Range("M96").Value = -632.09095
$ws.Range("N96").Value = -5240.5
$ws.Range("H99").Value = 92367.37
$ws.Range("J99").Value = 145012.42
$ws.Range("L99").Value = 435037.26
$ws.Range("N99").Value = -438033.26
$ws.Range("H101").Value = 742
$ws.Range("I101").Value = 767.4
$ws.Range("K101").Value = 2302.2
$ws.Range("M101").Value = -680.1999999999998
$ws.Range("H103").Value = 892.5714
$ws.Range("J103").Value = 900
$ws.Range("L103").Value = 2700
$ws.Range("N103").Value = -3872
$ws.Range("H106").Value = 3833.4092
$ws.Range("I106").Value = 3713.7058
$ws.Range("K106").Value = 3713.7058
$ws.Range("M106").Value = -3082.7058
$ws.Range("H112").Value = 2274.8333
$ws.Range("J112").Value = 2379.2942
$ws.Range("L112").Value = 7137.882599999999
$ws.Range("N112").Value = -9353.882599999999
$ws.Range("H118").Value = 527.5
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H125").Value = 11115040
$ws.Range("I125").Value = 763
$ws.Range("K125").Value = 6867
$ws.Range("M125").Value = -4407
$ws.Range("H132").Value = 2073.9583
$ws.Range("I132").Value = 1810.2727
$ws.Range("J132").Value = 4974.5
$ws.Range("K132").Value = 5430.8181
$ws.Range("L132").Value = 14923.5
$ws.Range("M132").Value = -2900.8181
$ws.Range("N132").Value = -19983.5
$ws.Range("H137").Value = 3311.238
$ws.Range("I137").Value = 2519.6897
$ws.Range("J137").Value = 5077
$ws.Range("K137").Value = 7559.0691
$ws.Range("L137").Value = 15231
$ws.Range("M137").Value = -5009.0691
$ws.Range("N137").Value = -20331
$ws.Range("H138").Value = 6593.68
$ws.Range("I138").Value = 3509.7058
$ws.Range("K138").Value = 10529.1174
$ws.Range("M138").Value = -5389.117400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1639
$ws.Range("I97").Value = 1639
$ws.Range("K97").Value = 1639
$ws.Range("M97").Value = -1143
$ws.Range("H102").Value = 1733.1212
$ws.Range("I102").Value = 1713.2069
$ws.Range("J102").Value = 1877.5
$ws.Range("K102").Value = 1713.2069
$ws.Range("L102").Value = 1877.5
$ws.Range("M102").Value = -91.20689999999991
$ws.Range("N102").Value = -5121.5
$ws.Range("H132").Value = 4033.8044
$ws.Range("I132").Value = 3575.139
$ws.Range("J132").Value = 5685
$ws.Range("K132").Value = 10725.417
$ws.Range("L132").Value = 17055
$ws.Range("M132").Value = -8195.417000000001
$ws.Range("N132").Value = -22115
$ws.Range("H138").Value = 69999.5
$ws.Range("J138").Value = 69999.5
$ws.Range("L138").Value = 69999.5
$ws.Range("N138").Value = -80279.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 333339970
$ws.Range("J20").Value = 9950
$ws.Range("L20").Value = 9950
$ws.Range("N20").Value = -10444
$ws.Range("H22").Value = 326.66666
$ws.Range("I22").Value = 324
$ws.Range("J22").Value = 340
$ws.Range("K22").Value = 324
$ws.Range("L22").Value = 340
$ws.Range("M22").Value = -151
$ws.Range("N22").Value = -686
$ws.Range("H99").Value = 985.1429000000001
$ws.Range("I99").Value = 919.2
$ws.Range("J99").Value = 1150
$ws.Range("K99").Value = 919.2
$ws.Range("L99").Value = 1150
$ws.Range("M99").Value = 578.8
$ws.Range("N99").Value = -4146
$ws.Range("H105").Value = 84858.086
$ws.Range("I105").Value = 112477.555
$ws.Range("J105").Value = 1999.6666
$ws.Range("K105").Value = 112477.555
$ws.Range("L105").Value = 1999.6666
$ws.Range("M105").Value = -110730.555
$ws.Range("N105").Value = -5493.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 729.7
$ws.Range("I105").Value = 766.3333
$ws.Range("K105").Value = 766.3333
$ws.Range("M105").Value = 980.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 494.625
$ws.Range("J44").Value = 801.3333
$ws.Range("L44").Value = 2403.9999
$ws.Range("N44").Value = -3199.9999
$ws.Range("H134").Value = 2254.087
$ws.Range("I134").Value = 1092.2
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 3276.6
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = 1793.4
$ws.Range("N134").Value = -40140
$ws.Range("H137").Value = 6546.579
$ws.Range("J137").Value = 7191.846
$ws.Range("L137").Value = 21575.538
$ws.Range("N137").Value = -31775.538

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3563.5173
$ws.Range("I122").Value = 3382.7727
$ws.Range("K122").Value = 10148.3181
$ws.Range("M122").Value = -7698.3181
$ws.Range("H126").Value = 3189.5293
$ws.Range("I126").Value = 2302.2
$ws.Range("J126").Value = 4457.143
$ws.Range("K126").Value = 6906.599999999999
$ws.Range("L126").Value = 13371.429
$ws.Range("M126").Value = -4436.599999999999
$ws.Range("N126").Value = -18311.429
$ws.Range("H132").Value = 151316.08
$ws.Range("I132").Value = 9856.909
$ws.Range("J132").Value = 669999.7
$ws.Range("K132").Value = 29570.727
$ws.Range("L132").Value = 2009999.1
$ws.Range("M132").Value = -27040.727
$ws.Range("N132").Value = -2015059.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6932.8335
$ws.Range("I7").Value = 4900
$ws.Range("J7").Value = 7949.25
$ws.Range("K7").Value = 4900
$ws.Range("L7").Value = 7949.25
$ws.Range("M7").Value = -4788
$ws.Range("N7").Value = -8173.25
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -205
$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 500
$ws.Range("K27").Value = 500
$ws.Range("M27").Value = -393
$ws.Range("H93").Value = 4928.304
$ws.Range("I93").Value = 5000.3335
$ws.Range("J93").Value = 4849.727
$ws.Range("K93").Value = 5000.3335
$ws.Range("L93").Value = 4849.727
$ws.Range("M93").Value = -3752.3335
$ws.Range("N93").Value = -7345.727
$ws.Range("H94").Value = 65000
$ws.Range("J94").Value = 65000
$ws.Range("L94").Value = 65000
$ws.Range("N94").Value = -66352
$ws.Range("H100").Value = 5298.8
$ws.Range("I100").Value = 2495
$ws.Range("J100").Value = 5999.75
$ws.Range("K100").Value = 2495
$ws.Range("L100").Value = 5999.75
$ws.Range("M100").Value = -1954
$ws.Range("N100").Value = -7081.75
$ws.Range("H126").Value = 6932.8335
$ws.Range("I126").Value = 4900
$ws.Range("J126").Value = 7949.25
$ws.Range("K126").Value = 14700
$ws.Range("L126").Value = 23847.75
$ws.Range("M126").Value = -12230
$ws.Range("N126").Value = -28787.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 44538.5
$ws.Range("J54").Value = 44538.5
$ws.Range("L54").Value = 44538.5
$ws.Range("N54").Value = -45578.5
$ws.Range("H81").Value = 11114.071
$ws.Range("I81").Value = 3199.625
$ws.Range("K81").Value = 6399.25
$ws.Range("M81").Value = -5338.25
$ws.Range("H84").Value = 11114.071
$ws.Range("I84").Value = 3199.625
$ws.Range("K84").Value = 31996.25
$ws.Range("M84").Value = -26692.25
$ws.Range("H126").Value = 1832.5834
$ws.Range("I126").Value = 1799.1
$ws.Range("K126").Value = 5397.299999999999
$ws.Range("M126").Value = -2927.299999999999
